$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'41.886.37"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +5.23%  "
$ws.Range("D3").Value = "'2.268.99"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.09%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("E5").Value = "  +3.68%  "
$ws.Range("D6").Value = "'92.85"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.69%  "
$ws.Range("E7").Value = "  +3.47%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'0.487"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.50%  "
$ws.Range("D10").Value = "'32.68"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.57%  "
$ws.Range("D11").Value = "'54.47"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +8.89%  "
$ws.Range("E12").Value = "  +2.98%  "
$ws.Range("E13").Value = "  +2.98%  "
$ws.Range("D14").Value = "'6.69"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.14%  "
$ws.Range("D15").Value = "'2.621.45"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.27%  "
$ws.Range("D16").Value = "'14.18"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.68%  "
$ws.Range("D17").Value = "'2.261.20"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.86%  "
$ws.Range("E18").Value = "  +3.81%  "
$ws.Range("D19").Value = "'41.802.76"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.28%  "
$ws.Range("D20").Value = "'12.32"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +10.58%  "
$ws.Range("E21").Value = "  +3.02%  "
$ws.Range("D22").Value = "'5.94"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.88%  "
$ws.Range("D23").Value = "'67.31"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.32%  "
$ws.Range("D24").Value = "'241.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.56%  "
$ws.Range("E25").Value = "  +5.71%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("E27").Value = "  +5.34%  "
$ws.Range("D28").Value = "'23.86"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.80%  "
$ws.Range("D29").Value = "'2.19"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.01%  "
$ws.Range("E30").Value = "  +5.75%  "
$ws.Range("D31").Value = "'34.24"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +9.59%  "
$ws.Range("D32").Value = "'158.26"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.48%  "
$ws.Range("D33").Value = "'1.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.10%  "
$ws.Range("D34").Value = "'5.18"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.29%  "
$ws.Range("E35").Value = "  +5.02%  "
$ws.Range("D36").Value = "'3.06"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.94%  "
$ws.Range("E37").Value = "  +3.20%  "
$ws.Range("D38").Value = "'0.105"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.70%  "
$ws.Range("D39").Value = "'16.51"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.82%  "
$ws.Range("E40").Value = "  +3.33%  "
$ws.Range("E41").Value = "  +6.71%  "
$ws.Range("D42").Value = "'3.98"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.94%  "
$ws.Range("D43").Value = "'20.41"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +17.46%  "
$ws.Range("D44").Value = "'2.050.48"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.90%  "
$ws.Range("D45").Value = "'0.0278"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.19%  "
$ws.Range("E46").Value = "  +3.23%  "
$ws.Range("E47").Value = "  +10.44%  "
$ws.Range("E48").Value = "  -4.19%  "
$ws.Range("D49").Value = "'2.493.89"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.50%  "
$ws.Range("E50").Value = "  +3.20%  "
$ws.Range("E51").Value = "  +4.57%  "
